$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-wrap the survey-question text in A4 (line break moved after "all countries")
$ws.Range("A4").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""

# Refreshed data values from the latest prepare/render run
$ws.Range("B2").Value = 0.697581258153402
$ws.Range("K2").Value = 0.686659119892822
$ws.Range("L2").Value = 0.745847844531965
$ws.Range("N2").Value = 0.617013940284116

$ws.Range("B3").Value = 0.641096347070675
$ws.Range("K3").Value = 0.552386382607627
$ws.Range("L3").Value = 0.752151504159109
$ws.Range("N3").Value = 0.575861430623479

$ws.Range("B4").Value = 0.717914385961719
$ws.Range("K4").Value = 0.702372413171302
$ws.Range("L4").Value = 0.770988593693527
$ws.Range("N4").Value = 0.561039368985046

$ws.Range("B5").Value = 0.456318549602673
$ws.Range("N5").Value = 0.402452424604714

$ws.Range("B6").Value = 0.604334051757766
$ws.Range("K6").Value = 0.584541124690159
$ws.Range("L6").Value = 0.568089649263453
$ws.Range("N6").Value = 0.55410364808233

$ws.Range("B7").Value = 0.55719177445442
